# Apply the crypto price/volume updates for Sun Jul 30 10:25:03 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.328.99'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.878.17'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''0.7286'
$ws.Range('E5').Value = '  +2.98%  '
$ws.Range('D6').Value = '''242.69'
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('D7').Value = '''1.002'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '''0.08003'
$ws.Range('E8').Value = '  +3.11%  '
$ws.Range('D9').Value = '''0.3162'
$ws.Range('E9').Value = '  +2.30%  '
$ws.Range('D10').Value = '''25.02'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').Value = '''0.08238'
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('D12').Value = '1.888.69'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').Value = '''94.65'
$ws.Range('E13').Value = '  +4.04%  '
$ws.Range('D14').Value = '''5.231'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').Value = '''0.7129'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').Value = '''6.399'
$ws.Range('E16').Value = '  +5.59%  '
$ws.Range('D17').Value = '''0.000008487'
$ws.Range('E17').Value = '  +3.85%  '
$ws.Range('D18').Value = '29.318.82'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').Value = '''243.66'
$ws.Range('E19').Value = '  +1.61%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '''13.25'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.129.85'
$ws.Range('E21').Value = '  +0.53%  '
$ws.Range('D22').Value = '''1.002'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').Value = '''7.757'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').Value = '''0.1613'
$ws.Range('E25').Value = '  +1.82%  '
$ws.Range('D26').Value = '''162.66'
$ws.Range('D27').Value = '''9.036'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').Value = '''18.55'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '''4.406'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').Value = '''4.307'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('D32').Value = '''1.187'
$ws.Range('E32').Value = '  -8.11%  '
$ws.Range('D33').Value = '''0.05374'
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('D34').Value = '''1.942'
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('D35').Value = '''0.7596'
$ws.Range('E35').Value = '  +2.18%  '
$ws.Range('D36').Value = '''1.178'
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('D38').Value = '''0.01878'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('D39').Value = '1.281.99'
$ws.Range('E39').Value = '  +4.10%  '
$ws.Range('D40').Value = '''2.761'
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('D41').Value = '''6.416'
$ws.Range('E41').Value = '  -2.22%  '
$ws.Range('D42').Value = '''113.29'
$ws.Range('E42').Value = '  +3.15%  '
$ws.Range('D43').Value = '''0.9071'
$ws.Range('E43').Value = '  +2.61%  '
$ws.Range('D44').Value = '''74.28'
$ws.Range('E44').Value = '  +2.62%  '
$ws.Range('E45').Value = '  +9.59%  '
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').Value = '2.024.85'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').Value = '''0.5228'
$ws.Range('E48').Value = '  +0.70%  '
$ws.Range('D49').Value = '''1.797'
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').Value = '''9.473'
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('D51').Value = '''0.4355'
$ws.Range('E51').Value = '  +1.16%  '
